$wb = $excel.ActiveWorkbook

# 1. Update "Generated on" timestamp on Summary sheet (B9)
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B9").Value = "2025-01-08 12:51:38"

# 2. Update the big dictionary-of-results string on Summary sheet (B20)
$summary.Range("B20").Value = '{''mcar'': {''mean'': {''bayesian'': {''0.05'': {''scores'': {''RMSE'': 0.5197922283008971, ''MAE'': 0.4543356516868202, ''MI'': 0.0, ''CORRELATION'': 0}, ''times'': {''contamination'': 0.00027060508728027344, ''optimization'': 0, ''imputation'': 0.0001704692840576172, ''log_imputation'': -8.676955429624789}}, ''0.1'': {''scores'': {''RMSE'': 1.0659202645786816, ''MAE'': 0.9085417731383956, ''MI'': 0.0, ''CORRELATION'': 0}, ''times'': {''contamination'': 0.00019931793212890625, ''optimization'': 0, ''imputation'': 0.00011801719665527344, ''log_imputation'': -9.044680209750107}}, ''0.2'': {''scores'': {''RMSE'': 1.1400385999631493, ''MAE'': 0.9394950730289477, ''MI'': 0.0, ''CORRELATION'': 0}, ''times'': {''contamination'': 0.00022220611572265625, ''optimization'': 0, ''imputation'': 9.34600830078125e-05, ''log_imputation'': -9.277976132528334}}, ''0.4'': {''scores'': {''RMSE'': 1.0333061850175014, ''MAE'': 0.817160720129779, ''MI'': 0.0, ''CORRELATION'': 0}, ''times'': {''contamination'': 0.00030231475830078125, ''optimization'': 0, ''imputation'': 6.556510925292969e-05, ''log_imputation'': -9.632466874652225}}, ''0.6'': {''scores'': {''RMSE'': 1.0938413270459857, ''MAE'': 0.8545290213993658, ''MI'': 0.0, ''CORRELATION'': 0}, ''times'': {''contamination'': 0.0005660057067871094, ''optimization'': 0, ''imputation'': 6.151199340820312e-05, ''log_imputation'': -9.696278387397179}}, ''0.8'': {''scores'': {''RMSE'': 1.07436956341757, ''MAE'': 0.8291370178635111, ''MI'': 0.0, ''CORRELATION'': 0}, ''times'': {''contamination'': 0.0012047290802001953, ''optimization'': 0, ''imputation'': 6.29425048828125e-05, ''log_imputation'': -9.673288869172481}}}}, ''cdrec'': {''bayesian'': {''0.05'': {''scores'': {''RMSE'': 0.37483452324301586, ''MAE'': 0.3375262694281006, ''MI'': 1.0397207708399179, ''CORRELATION'': 0.7365655689896633}, ''times'': {''contamination'': 0.00012755393981933594, ''optimization'': 0.5108654499053955, ''imputation'': 0.004567623138427734, ''log_imputation'': -5.388762310375959}}, ''0.1'': {''scores'': {''RMSE'': 1.3799678230195285, ''MAE'': 1.1003322284844623, ''MI'': 1.732867951399863, ''CORRELATION'': -0.500100644242659}, ''times'': {''contamination'': 0.00033164024353027344, ''optimization'': 0.5108654499053955, ''imputation'': 0.004395723342895508, ''log_imputation'': -5.427123178248406}}, ''0.2'': {''scores'': {''RMSE'': 0.5279485898506157, ''MAE'': 0.42431581904234256, ''MI'': 1.342409426595628, ''CORRELATION'': 0.9071070625126642}, ''times'': {''contamination'': 0.00019049644470214844, ''optimization'': 0.5108654499053955, ''imputation'': 0.011562585830688477, ''log_imputation'': -4.459980752976466}}, ''0.4'': {''scores'': {''RMSE'': 0.6529812630837011, ''MAE'': 0.42858056477338186, ''MI'': 0.8905639332827393, ''CORRELATION'': 0.7722811146383882}, ''times'': {''contamination'': 0.00030994415283203125, ''optimization'': 0.5108654499053955, ''imputation'': 0.00451207160949707, ''log_imputation'': -5.4009988940131475}}, ''0.6'': {''scores'': {''RMSE'': 0.6798826891423311, ''MAE'': 0.47135122101632737, ''MI'': 0.6001676421795947, ''CORRELATION'': 0.7742382236368857}, ''times'': {''contamination'': 0.0005815029144287109, ''optimization'': 0.5108654499053955, ''imputation'': 0.0029671192169189453, ''log_imputation'': -5.820163757517104}}, ''0.8'': {''scores'': {''RMSE'': 0.7608485588056992, ''MAE'': 0.5479154581689161, ''MI'': 0.42721564894947844, ''CORRELATION'': 0.7017141157422242}, ''times'': {''contamination'': 0.0012214183807373047, ''optimization'': 0.5108654499053955, ''imputation'': 0.004498004913330078, ''log_imputation'': -5.404121333108998}}}}, ''stmvl'': {''bayesian'': {''0.05'': {''scores'': {''RMSE'': 0.3251125774837754, ''MAE'': 0.26797673641099284, ''MI'': 1.0397207708399179, ''CORRELATION'': 0.6142581896031455}, ''times'': {''contamination'': 0.00016355514526367188, ''optimization'': 0.18949556350708008, ''imputation'': 0.0333857536315918, ''log_imputation'': -3.399626007968974}}, ''0.1'': {''scores'': {''RMSE'': 0.299492451492057, ''MAE'': 0.26432871720074347, ''MI'': 1.9061547465398494, ''CORRELATION'': 0.967896575643492}, ''times'': {''contamination'': 0.00013828277587890625, ''optimization'': 0.18949556350708008, ''imputation'': 0.034218549728393555, ''log_imputation'': -3.374987392234605}}, ''0.2'': {''scores'': {''RMSE'': 0.32852543256899075, ''MAE'': 0.27202573018354975, ''MI'': 1.5996631161656454, ''CORRELATION'': 0.9558373872353643}, ''times'': {''contamination'': 0.00016021728515625, ''optimization'': 0.18949556350708008, ''imputation'': 0.03754568099975586, ''log_imputation'': -3.282196927367169}}, ''0.4'': {''scores'': {''RMSE'': 0.4508488005700101, ''MAE'': 0.34941433537269606, ''MI'': 0.8543113555966528, ''CORRELATION'': 0.8959297471926679}, ''times'': {''contamination'': 0.0002884864807128906, ''optimization'': 0.18949556350708008, ''imputation'': 0.03878617286682129, ''log_imputation'': -3.2496914652794007}}, ''0.6'': {''scores'': {''RMSE'': 18.797539991079297, ''MAE'': 7.812583796335101, ''MI'': 0.36244773022350796, ''CORRELATION'': 0.6210142190959098}, ''times'': {''contamination'': 0.0006196498870849609, ''optimization'': 0.18949556350708008, ''imputation'': 0.029697656631469727, ''log_imputation'': -3.51668713758367}}, ''0.8'': {''scores'': {''RMSE'': 3.1451455567216193, ''MAE'': 1.1637520656636082, ''MI'': 0.0643204354315137, ''CORRELATION'': 0.22737088719870605}, ''times'': {''contamination'': 0.0011970996856689453, ''optimization'': 0.18949556350708008, ''imputation'': 0.034436941146850586, ''log_imputation'': -3.368625420241316}}}}, ''iim'': {''bayesian'': {''0.05'': {''scores'': {''RMSE'': 0.2311363556202525, ''MAE'': 0.22809317150257158, ''MI'': 0.6931471805599452, ''CORRELATION'': 0.8754093900930757}, ''times'': {''contamination'': 0.00020766258239746094, ''optimization'': 0.8132822513580322, ''imputation'': 0.006684064865112305, ''log_imputation'': -5.008028963761024}}, ''0.1'': {''scores'': {''RMSE'': 0.21734571962767568, ''MAE'': 0.20142183555276616, ''MI'': 1.4941751382893083, ''CORRELATION'': 0.9836625389334559}, ''times'': {''contamination'': 0.0001220703125, ''optimization'': 0.8132822513580322, ''imputation'': 0.008673667907714844, ''log_imputation'': -4.747463520246684}}, ''0.2'': {''scores'': {''RMSE'': 0.2763681623559098, ''MAE'': 0.21205899863451294, ''MI'': 1.692828654044598, ''CORRELATION'': 0.9663556239228223}, ''times'': {''contamination'': 0.00015807151794433594, ''optimization'': 0.8132822513580322, ''imputation'': 0.02692556381225586, ''log_imputation'': -3.6146791160565725}}, ''0.4'': {''scores'': {''RMSE'': 0.32470532661816204, ''MAE'': 0.24836184775095202, ''MI'': 1.0631520030142667, ''CORRELATION'': 0.9435024215665483}, ''times'': {''contamination'': 0.00028705596923828125, ''optimization'': 0.8132822513580322, ''imputation'': 0.19279956817626953, ''log_imputation'': -1.6461041365578395}}, ''0.6'': {''scores'': {''RMSE'': 0.45693859713260937, ''MAE'': 0.3350566242376081, ''MI'': 0.836724518636222, ''CORRELATION'': 0.9015668975756113}, ''times'': {''contamination'': 0.0006051063537597656, ''optimization'': 0.8132822513580322, ''imputation'': 0.5096421241760254, ''log_imputation'': -0.6740465168845}}, ''0.8'': {''scores'': {''RMSE'': 0.7301676007328138, ''MAE'': 0.5391664379693699, ''MI'': 0.43198783605819785, ''CORRELATION'': 0.7329833767632488}, ''times'': {''contamination'': 0.001192331314086914, ''optimization'': 0.8132822513580322, ''imputation'': 2.1703436374664307, ''log_imputation'': 0.7748855132944805}}}}, ''mrnn'': {''bayesian'': {''0.05'': {''scores'': {''RMSE'': 1.5534821982949707, ''MAE'': 1.5060357479863002, ''MI'': 0.0, ''CORRELATION'': 0}, ''times'': {''contamination'': 0.00022554397583007812, ''optimization'': 33.18070316314697, ''imputation'': 7.837775945663452, ''log_imputation'': 2.0589551137173507}}, ''0.1'': {''scores'': {''RMSE'': 1.2542309808619898, ''MAE'': 1.077694959981515, ''MI'': 0.0, ''CORRELATION'': 0}, ''times'': {''contamination'': 0.0003402233123779297, ''optimization'': 33.18070316314697, ''imputation'': 7.026872158050537, ''log_imputation'': 1.9497416790817856}}, ''0.2'': {''scores'': {''RMSE'': 1.5269073369222408, ''MAE'': 1.229707427832888, ''MI'': 0.397640786553919, ''CORRELATION'': 0.0010088946702063234}, ''times'': {''contamination'': 0.00017547607421875, ''optimization'': 33.18070316314697, ''imputation'': 6.481844186782837, ''log_imputation'': 1.8690050665999636}}, ''0.4'': {''scores'': {''RMSE'': 1.3382262981163862, ''MAE'': 1.079772579102801, ''MI'': 0.16557034901046197, ''CORRELATION'': 0.016391846298721195}, ''times'': {''contamination'': 0.0003674030303955078, ''optimization'': 33.18070316314697, ''imputation'': 8.01118803024292, ''log_imputation'': 2.080839068464571}}, ''0.6'': {''scores'': {''RMSE'': 1.3571781743580735, ''MAE'': 1.0989570438570895, ''MI'': 0.10967079612419346, ''CORRELATION'': -0.015125969685014458}, ''times'': {''contamination'': 0.0014710426330566406, ''optimization'': 33.18070316314697, ''imputation'': 9.495571374893188, ''log_imputation'': 2.2508255188511095}}, ''0.8'': {''scores'': {''RMSE'': 1.1992433951193446, ''MAE'': 0.9771615695053146, ''MI'': 0.050276623880035565, ''CORRELATION'': 0.00024685644553336644}, ''times'': {''contamination'': 0.0013470649719238281, ''optimization'': 33.18070316314697, ''imputation'': 8.68370532989502, ''log_imputation'': 2.1614483189086586}}}}}}'

# 3. RMSE sheet - mrnn row (F34:F39)
$rmse = $wb.Worksheets.Item("RMSE")
$rmseCells = @("F34","F35","F36","F37","F38","F39")
$rmseValues = @("1.5534821983","1.2542309809","1.5269073369","1.3382262981","1.3571781744","1.1992433951")
for ($i = 0; $i -lt $rmseCells.Length; $i++) {
    $cell = $rmse.Range($rmseCells[$i])
    $cell.NumberFormat = "@"
    $cell.Value = $rmseValues[$i]
    $cell.Style = "Normal"
}

# 4. MAE sheet - mrnn row (F34:F39)
$mae = $wb.Worksheets.Item("MAE")
$maeCells = @("F34","F35","F36","F37","F38","F39")
$maeValues = @("1.5060357480","1.0776949600","1.2297074278","1.0797725791","1.0989570439","0.9771615695")
for ($i = 0; $i -lt $maeCells.Length; $i++) {
    $cell = $mae.Range($maeCells[$i])
    $cell.NumberFormat = "@"
    $cell.Value = $maeValues[$i]
    $cell.Style = "Normal"
}

# 5. MI sheet - mrnn row, only the 0.4/0.6/0.8 thresholds changed (F37:F39)
$mi = $wb.Worksheets.Item("MI")
$miCells = @("F37","F38","F39")
$miValues = @("0.1655703490","0.1096707961","0.0502766239")
for ($i = 0; $i -lt $miCells.Length; $i++) {
    $cell = $mi.Range($miCells[$i])
    $cell.NumberFormat = "@"
    $cell.Value = $miValues[$i]
    $cell.Style = "Normal"
}

# 6. CORRELATION sheet - mrnn row, only the 0.2/0.4/0.6/0.8 thresholds changed (F36:F39)
$correlation = $wb.Worksheets.Item("CORRELATION")
$correlationCells = @("F36","F37","F38","F39")
$correlationValues = @("0.0010088947","0.0163918463","-0.0151259697","0.0002468564")
for ($i = 0; $i -lt $correlationCells.Length; $i++) {
    $cell = $correlation.Range($correlationCells[$i])
    $cell.NumberFormat = "@"
    $cell.Value = $correlationValues[$i]
    $cell.Style = "Normal"
}
